$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 48587
$ws.Range("C2").Value = 48587
$ws.Range("D2").Value = 48587
$ws.Range("E2").Value = 48587
$ws.Range("F2").Value = 48587
$ws.Range("G2").Value = 48587
$ws.Range("H2").Value = 48587
$ws.Range("I2").Value = 48587
$ws.Range("J2").Value = 48587
$ws.Range("K2").Value = 48587
$ws.Range("L2").Value = 48587
$ws.Range("M2").Value = 48587
$ws.Range("N2").Value = 48587
$ws.Range("O2").Value = 48587
$ws.Range("P2").Value = 48587
$ws.Range("Q2").Value = 48587

$ws.Range("B3").Value = 18482.26115607056
$ws.Range("C3").Value = 2.945664478152592
$ws.Range("D3").Value = 43.23594788729496
$ws.Range("E3").Value = 0.07080083149813736
$ws.Range("F3").Value = 0.3693580587399922
$ws.Range("G3").Value = 130943.1803568856
$ws.Range("H3").Value = 67.66391483318583
$ws.Range("I3").Value = 3.217689361351802
$ws.Range("J3").Value = 2507.886269372466
$ws.Range("K3").Value = 913.8347380986684
$ws.Range("L3").Value = 462.4929754872703
$ws.Range("M3").Value = 1157.553899602774
$ws.Range("N3").Value = 1631.480181941672
$ws.Range("O3").Value = 198.3287240922613
$ws.Range("P3").Value = 6.949048243357278
$ws.Range("Q3").Value = 6.450532035318089

$ws.Range("B4").Value = 28191.71108331252
$ws.Range("C4").Value = 1.411945396998974
$ws.Range("D4").Value = 29.69981647724367
$ws.Range("E4").Value = 0.2564944985987722
$ws.Range("F4").Value = 0.4826359678043436
$ws.Range("G4").Value = 75208.80070856874
$ws.Range("H4").Value = 14.93884860387811
$ws.Range("I4").Value = 0.4233053916825096
$ws.Range("J4").Value = 5739.409853448393
$ws.Range("K4").Value = 5662.93151213163
$ws.Range("L4").Value = 5213.825652896686
$ws.Range("M4").Value = 3860.971557146052
$ws.Range("N4").Value = 3466.857487148857
$ws.Range("O4").Value = 36.14278185225009
$ws.Range("P4").Value = 1.03435187816563
$ws.Range("Q4").Value = 3.250982181878554

$ws.Range("H5").Value = 28.84
$ws.Range("I5").Value = 2.514

$ws.Range("B6").Value = 2207.065
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 19
$ws.Range("G6").Value = 37392
$ws.Range("H6").Value = 56.36
$ws.Range("I6").Value = 2.759
$ws.Range("O6").Value = 210.8967606
$ws.Range("P6").Value = 6.489

$ws.Range("B7").Value = 7358.74
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 36
$ws.Range("G7").Value = 151315
$ws.Range("H7").Value = 69.36
$ws.Range("I7").Value = 3.29
$ws.Range("O7").Value = 214.7027646
$ws.Range("P7").Value = 7.143

$ws.Range("B8").Value = 21254.705
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 71
$ws.Range("G8").Value = 202307
$ws.Range("H8").Value = 80.84
$ws.Range("I8").Value = 3.594
$ws.Range("J8").Value = 2797.26
$ws.Range("K8").Value = 7.64
$ws.Range("L8").Value = 4.54
$ws.Range("M8").Value = 536.8200000000001
$ws.Range("N8").Value = 2093.48
$ws.Range("O8").Value = 219.8118854
$ws.Range("P8").Value = 7.808

$ws.Range("B9").Value = 385051.04
$ws.Range("C9").Value = 5
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 205863
$ws.Range("H9").Value = 93.34
$ws.Range("I9").Value = 3.907
$ws.Range("J9").Value = 75149.78999999999
$ws.Range("K9").Value = 92523.94
$ws.Range("L9").Value = 83340.33
$ws.Range("M9").Value = 48159.86
$ws.Range("N9").Value = 36430.33
$ws.Range("O9").Value = 226.9873637
$ws.Range("P9").Value = 8.622999999999999
